$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Content changes -----------------------------------------------------
# Header label "CSS Id" -> "VLJ #"
$ws.Range("C2").Value = "VLJ #"

# Example judge id placeholder "BVAJONESB" -> "123" (repeated down the
# example rows, same as "Jones, Bernard" repeats in column B)
$ws.Range("C3:C7").Value = "123"

# --- Footer font: "Helvetica Neue" -> "Helvetica" -------------------------
$ws.PageSetup.CenterFooter = "&""Helvetica,Regular""&12&K000000&P"

# --- Workbook theme font: "Helvetica Neue" -> "Helvetica" -----------------
$themeFonts = $wb.Theme.ThemeFontScheme
$themeFonts.MajorFont.Latin = "Helvetica"
$themeFonts.MinorFont.Latin = "Helvetica"

# --- Two new footer rows (9 and 10) ---------------------------------------
# These extend the table with a thin red outline box below row 8, matching
# the look of the red highlight borders used elsewhere in the template.
$red = 255

$ws.Range("A9:I10").RowHeight = 17

$row9 = $ws.Range("A9:I9")
$row9.Interior.ColorIndex = 0
$row9.Font.Name = "Calibri"
$row9.Font.Size = 12

$row10 = $ws.Range("A10:I10")
$row10.Interior.ColorIndex = 0
$row10.Font.Name = "Calibri"
$row10.Font.Size = 12

# Row 9 - top edge of the red outline box
$a9 = $ws.Range("A9")
$a9.Borders.Item(7).Color = $red
$a9.Borders.Item(7).LineStyle = 1
$a9.Borders.Item(8).Color = $red
$a9.Borders.Item(8).LineStyle = 1

$bh9 = $ws.Range("B9:H9")
$bh9.Borders.Item(8).Color = $red
$bh9.Borders.Item(8).LineStyle = 1

$i9 = $ws.Range("I9")
$i9.Borders.Item(8).Color = $red
$i9.Borders.Item(8).LineStyle = 1
$i9.Borders.Item(10).Color = $red
$i9.Borders.Item(10).LineStyle = 1

# Row 10 - bottom edge of the red outline box
$a10 = $ws.Range("A10")
$a10.Borders.Item(7).Color = $red
$a10.Borders.Item(7).LineStyle = 1
$a10.Borders.Item(9).Color = $red
$a10.Borders.Item(9).LineStyle = 1

$bh10 = $ws.Range("B10:H10")
$bh10.Borders.Item(9).Color = $red
$bh10.Borders.Item(9).LineStyle = 1

$i10 = $ws.Range("I10")
$i10.Borders.Item(9).Color = $red
$i10.Borders.Item(9).LineStyle = 1
$i10.Borders.Item(10).Color = $red
$i10.Borders.Item(10).LineStyle = 1
